# Applies crypto price/volume/coin-order update
# Commit message: Updated symbol list on Sat Feb 11 17:15:00 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cellUpdates = @(
    @('D2', '309.59'),
    @('E2', '0.73%'),
    @('D3', '40.81'),
    @('E3', '0.99%'),
    @('D4', '5.126'),
    @('E4', '1.73%'),
    @('D5', '0.07630'),
    @('E5', '0.31%'),
    @('B6', 'FTXToken'),
    @('C6', 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'),
    @('D6', '1.623'),
    @('E6', '1.72%'),
    @('B7', 'BTSEToken'),
    @('C7', 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'),
    @('D7', '2.499'),
    @('E7', '2.24%'),
    @('B8', 'MXToken'),
    @('C8', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'),
    @('D8', '0.9096'),
    @('E8', '0.19%'),
    @('B9', 'LiechtensteinCryptoassetsExchange'),
    @('C9', 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'),
    @('D9', '0.1194'),
    @('E9', '19.32%'),
    @('B10', 'WazirX'),
    @('C10', 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'),
    @('D10', '0.1820'),
    @('E10', '4.10%'),
    @('B11', 'MandalaExchangeToken'),
    @('C11', 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'),
    @('D11', '0.09147'),
    @('E11', '1.39%'),
    @('B12', 'BitrueCoin'),
    @('C12', 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'),
    @('D12', '0.04249'),
    @('E12', '-2.03%'),
    @('B13', 'BitMartToken'),
    @('C13', 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'),
    @('D13', '0.1046'),
    @('E13', '-0.82%'),
    @('B14', 'BitForexToken'),
    @('C14', 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'),
    @('D14', '0.001252'),
    @('E14', '1.92%'),
    @('B15', 'TigerCash'),
    @('C15', 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'),
    @('D15', '0.005797'),
    @('E15', '-0.99%'),
    @('B16', 'LEO'),
    @('C16', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @('D16', '3.356'),
    @('E16', '-0.42%'),
    @('B17', 'GateToken'),
    @('C17', 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'),
    @('D17', '4.283'),
    @('E17', '0.95%'),
    @('E18', '-0.66%'),
    @('D19', '6.907'),
    @('E19', '1.18%'),
    @('E20', '4.36%'),
    @('D21', '0.2707'),
    @('E21', '-5.02%'),
    @('D22', '0.04052'),
    @('E22', '-2.50%'),
    @('E23', '4.46%'),
    @('D24', '0.004105'),
    @('E24', '1.04%'),
    @('D25', '0.0001273'),
    @('E25', '-2.29%'),
    @('D26', '0.0003752'),
    @('D38', '0.02432'),
    @('E38', '0.49%'),
    @('D39', '0.05237'),
    @('E39', '2.04%'),
    @('D40', '0.007787'),
    @('E40', '-0.86%'),
    @('D41', '0.1301'),
    @('E41', '-0.23%'),
    @('D42', '0.006807'),
    @('E42', '-4.24%'),
    @('D43', '0.001934'),
    @('E43', '-0.92%'),
    @('D44', '0.008081'),
    @('E44', '-3.42%'),
    @('D45', '0.3069'),
    @('E45', '-7.49%'),
    @('D46', '0.00006902'),
    @('E46', '7.00%'),
    @('D47', '0.00000000752'),
    @('E47', '0.01%'),
    @('D48', '0.09607'),
    @('E48', '1,700.35%'),
    @('D50', '0.00002105'),
    @('E50', '0.01%'),
    @('D51', '0.0002004'),
    @('E51', '0.01%')
)

foreach ($pair in $cellUpdates) {
    $cellRef = $pair[0]
    $text = $pair[1]
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

